$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as text so
# values like "321.30" or "2.40" keep their exact literal formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '48.311.48'
$ws.Range("E2").Value = '  +2.03%  '
$ws.Range("D3").Value = '2.506.04'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("D5").Value = '321.30'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '108.06'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("D7").Value = '0.529'
$ws.Range("E7").Value = '  +1.35%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").Value = '39.94'
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("D11").Value = '20.28'
$ws.Range("E11").Value = '  +9.20%  '
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = '2.898.27'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").Value = '2.508.41'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '0.846'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '48.155.02'
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("D19").Value = '13.11'
$ws.Range("E19").Value = '  -2.42%  '
$ws.Range("D20").Value = '6.74'
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '282.59'
$ws.Range("E23").Value = '  +14.45%  '
$ws.Range("D24").Value = '72.33'
$ws.Range("E24").Value = '  +2.39%  '
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '25.74'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '2.40'
$ws.Range("E28").Value = '  +4.72%  '
$ws.Range("D29").Value = '9.79'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("D31").Value = '35.28'
$ws.Range("E31").Value = '  +2.11%  '
$ws.Range("D32").Value = '49.36'
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("D33").Value = '19.61'
$ws.Range("E33").Value = '  -3.79%  '
$ws.Range("D34").Value = '5.36'
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '0.0783'
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '4.65'
$ws.Range("E38").Value = '  -2.08%  '
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").Value = '121.52'
$ws.Range("E41").Value = '  +2.07%  '
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("D43").Value = '21.54'
$ws.Range("E43").Value = '  -4.38%  '
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("D45").Value = '2.014.81'
$ws.Range("E45").Value = '  +1.14%  '
$ws.Range("D46").Value = '3.17'
$ws.Range("E46").Value = '  +4.66%  '
$ws.Range("D47").Value = '1.85'
$ws.Range("E47").Value = '  +3.89%  '
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("D49").Value = '8.99'
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("D50").Value = '5.18'
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").Value = '80.48'
$ws.Range("E51").Value = '  +3.86%  '
